# Sprint Backlog.xlsx update
# - A2 used to hold a broken dynamic-array formula "=- Martim Costa" (#NAME? error).
#   Replace it with the plain text "Martim Costa".
# - A4's label loses its leading " - " separator: "Identifie 3 code smells - Martim Costa"
#   becomes "Identifie 3 code smells Martim Costa".
# - D4 (previously blank) now carries the label "Identifie 3 code smells".
# - Move the active selection to D20 (matches the saved sheetView selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Martim Costa"
$ws.Range("A4").Value = "Identifie 3 code smells Martim Costa"
$ws.Range("D4").Value = "Identifie 3 code smells"

$ws.Range("D20").Select()
